$wb = $excel.ActiveWorkbook

# --- Summary sheet ---
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 2      # Total IN
$summary.Range("B4").Value = 4      # Total OUT
$summary.Range("A10").Value = "Female"
$summary.Range("B10").Value = 1
$summary.Range("C10").Value = 1
$summary.Range("B16").Value = 1
$summary.Range("C16").Value = 2

# --- Hourly Breakdown sheet ---
$hourly = $wb.Worksheets.Item("Hourly Breakdown")
$hourly.Range("B18").Value = 2
$hourly.Range("C18").Value = 4
$hourly.Range("D18").Value = -2

# --- Charts Data sheet ---
$charts = $wb.Worksheets.Item("Charts Data")

# Insert a new row before row 4 (pushes Age Distribution block and below down by 1)
$charts.Rows.Item(4).Insert()
$charts.Range("A4").Value = "Female"
$charts.Range("B4").Value = 1

# Insert another new row before what is now row 9 (the row holding "26-45")
# After the first insert, original row6 (Age Distribution) is now row7,
# row7 (Age Group/Count headers) is now row8, row8 (26-45/1) is now row9.
$charts.Rows.Item(9).Insert()
$charts.Rows.Item(9).ClearFormats()
$charts.Range("A9").Value = "13-25"
$charts.Range("B9").Value = 1
